# Fruta / hortaliza, semanal
# The data rows (2-24) have been reordered (e.g. re-sorted by date); every
# cell in the "after" state is a straight copy of some other row's full
# contents (columns A:T) from the "before" state. Capture a snapshot of all
# rows first, then write the rows back out according to the permutation so
# that reads never see already-overwritten data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (both referring to the ORIGINAL/before layout)
$perm = @{
  2  = 15
  3  = 9
  4  = 5
  5  = 2
  6  = 3
  7  = 18
  8  = 19
  9  = 14
  10 = 8
  11 = 22
  12 = 23
  13 = 24
  14 = 20
  15 = 21
  16 = 17
  17 = 16
  18 = 11
  19 = 12
  20 = 13
  21 = 4
  22 = 10
  23 = 6
  24 = 7
}

$firstCol = 1   # A
$lastCol  = 20  # T

# Snapshot every data row (2-24) across columns A:T before any writes happen.
$snapshot = @{}
for ($r = 2; $r -le 24; $r++) {
  $rowVals = @()
  for ($c = $firstCol; $c -le $lastCol; $c++) {
    $rowVals += , ($ws.Cells.Item($r, $c).Value2)
  }
  $snapshot[$r] = $rowVals
}

# Write each destination row using the snapshot of its source row.
for ($r = 2; $r -le 24; $r++) {
  $srcRow = $perm[$r]
  $rowVals = $snapshot[$srcRow]
  for ($c = $firstCol; $c -le $lastCol; $c++) {
    $ws.Cells.Item($r, $c).Value2 = $rowVals[$c - $firstCol]
  }
}
